$wb = $excel.ActiveWorkbook

# --- Sheet "Input" (sheet1) ---
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("D3").Value = "Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat"
$wsInput.Columns.Item(4).ColumnWidth = 66.7109375

# --- Sheet "Calc" (sheet2) ---
$wsCalc = $wb.Worksheets.Item("Calc")

# New columns BG (59) and BH (60) with headers and meta row
$wsCalc.Cells.Item(1, 59).Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Cells.Item(1, 60).Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Cells.Item(2, 59).Value = "(%)"
$wsCalc.Cells.Item(2, 60).Value = "(%)"
$wsCalc.Cells.Item(3, 59).Value = 100
$wsCalc.Cells.Item(3, 60).Value = 100

# Copy header/meta cell styles from neighboring column BF (58) to new columns
$wsCalc.Cells.Item(1, 58).Copy()
$wsCalc.Cells.Item(1, 59).PasteSpecial(-4122) | Out-Null
$wsCalc.Cells.Item(1, 60).PasteSpecial(-4122) | Out-Null
$wsCalc.Cells.Item(2, 58).Copy()
$wsCalc.Cells.Item(2, 59).PasteSpecial(-4122) | Out-Null
$wsCalc.Cells.Item(2, 60).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Updated recalculated values in row 3
$wsCalc.Range("AP3").Value = 0.536
$wsCalc.Range("AQ3").Value = 0.1897992129706516
$wsCalc.Range("AW3").Value = 0.5304
$wsCalc.Range("AX3").Value = 0.5342210465052091
$wsCalc.Range("AY3").Value = 0.1878513580087756
$wsCalc.Range("BC3").Value = 0.5688527934319699
$wsCalc.Range("BE3").Value = 267.1105232526045
$wsCalc.Range("BF3").Value = 0.189204655095914

# Column width adjustments
$wsCalc.Columns.Item(43).ColumnWidth = 19.7109375
$wsCalc.Columns.Item(49).ColumnWidth = 9.7109375
$wsCalc.Columns.Item(51).ColumnWidth = 19.7109375
$wsCalc.Columns.Item(58).ColumnWidth = 18.7109375
$wsCalc.Columns.Item(59).ColumnWidth = 32.7109375
$wsCalc.Columns.Item(60).ColumnWidth = 30.7109375

# --- Sheet "Results" (sheet3) ---
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("N3").Value = 0.536
$wsResults.Range("P3").Value = 0.5304
$wsResults.Range("R3").Value = 0.5688527934319699
$wsResults.Columns.Item(16).ColumnWidth = 8.7109375

# --- Sheet "Constants" (sheet4) ---
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B3").Value = 0.00005
